$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tables")

$range = $ws.Range("A65:E94")
$key1 = $ws.Range("A65:A94")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 1, 0, 0)
$ws.Sort.SetRange($range)
$ws.Sort.Header = 0
$ws.Sort.Apply()
